$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 39)
$ws.Range("D2").Value = [double]"0.9999999999845814"
$ws.Range("E2").Value = [double]"0.9999999999845814"

# Row 3 (Control 17)
$ws.Range("D3").Value = [double]"0.004728029049891378"
$ws.Range("E3").Value = [double]"0.004728029049891378"

# Row 4 (Control 23)
$ws.Range("D4").Value = [double]"0.8472220674138857"
$ws.Range("E4").Value = [double]"0.8472220674138857"

# Row 5 (Control 27)
$ws.Range("D5").Value = [double]"5.112212561942627E-05"
$ws.Range("E5").Value = [double]"5.112212561942627E-05"

# Row 6 (Control 8)
$ws.Range("D6").Value = [double]"9.955964823549761E-15"
$ws.Range("E6").Value = [double]"9.955964823549761E-15"

# Row 7 (MDD 38)
$ws.Range("D7").Value = [double]"0.9978482232743362"
$ws.Range("E7").Value = [double]"0.002151776725663801"

# Row 8 (MDD 9)
$ws.Range("D8").Value = [double]"0.999934613030725"
$ws.Range("E8").Value = [double]"6.538696927504972E-05"

# Row 9 (MDD 49)
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = [double]"0.7885580150999533"
$ws.Range("E9").Value = [double]"0.2114419849000467"

# Row 10 (MDD 26)
$ws.Range("D10").Value = [double]"0.9999999999874902"
$ws.Range("E10").Value = [double]"1.250977099687134E-11"

# Row 11 (MDD 34)
$ws.Range("D11").Value = [double]"0.0003226301852358588"
$ws.Range("E11").Value = [double]"0.9996773698147642"
$ws.Range("F11").Value = [double]"3.505778074264526"
$ws.Range("G11").Value = [double]"0.7"
